# This edit inserts a brand new data row at row 159 ("semanal" weekly update
# of the Betarraga hortaliza price series for Vega Monumental Concepción).
# Inserting the row shifts all the existing rows 159-250 down to 160-251,
# which is exactly the change the commit's XML diff shows (every existing
# row's content simply moves one row down, and a new row 159 appears with
# fresh data while the sheet's used range grows from A1:R250 to A1:R251).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 159; this pushes the old
# rows 159..250 down to 160..251 and carries their formatting with them.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new record's values.
$ws.Cells.Item(159, 1).Value = 11
$ws.Cells.Item(159, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(159, 3).Value = "Bíobío"
$ws.Cells.Item(159, 4).Value = 44572
$ws.Cells.Item(159, 5).Value = 8
$ws.Cells.Item(159, 6).Value = 100114014
$ws.Cells.Item(159, 7).Value = "Betarraga"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 280
$ws.Cells.Item(159, 11).Value = 600
$ws.Cells.Item(159, 12).Value = 650
$ws.Cells.Item(159, 13).Value = 627
$ws.Cells.Item(159, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(159, 15).Value = "Región Metropolitana"
$ws.Cells.Item(159, 16).Value = 125
$ws.Cells.Item(159, 17).Value = 5
$ws.Cells.Item(159, 18).Value = "Hortaliza"
